$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = 0.03528612143816936
$ws.Range("I1").Value = 0.03528612143816936
$ws.Range("M1").Value = 0.03528612143816936
$ws.Range("O1").Value = 0.03528612143816936
$ws.Range("P1").Value = 0.03528612143816936
$ws.Range("R1").Value = 0.03528612143816936
$ws.Range("S1").Value = 0.03528612143816936
$ws.Range("U1").Value = 0.03528612143816936
$ws.Range("AK1").Value = 0.03528612143816936
$ws.Range("AM1").Value = 0.03528612143816936

$ws.Range("G2").Value = 0.006457780777315358
$ws.Range("I2").Value = 0.006457780777315358
$ws.Range("M2").Value = 0.006457780777315358
$ws.Range("O2").Value = 0.006457780777315358
$ws.Range("P2").Value = 0.006457780777315358
$ws.Range("R2").Value = 0.006457780777315358
$ws.Range("S2").Value = 0.006457780777315358
$ws.Range("U2").Value = 0.006457780777315358
$ws.Range("AK2").Value = 0.006457780777315358
$ws.Range("AM2").Value = 0.006457780777315358

$ws.Range("G3").Value = 0.00584934642747072
$ws.Range("I3").Value = 0.00584934642747072
$ws.Range("M3").Value = 0.00584934642747072
$ws.Range("O3").Value = 0.00584934642747072
$ws.Range("P3").Value = 0.00584934642747072
$ws.Range("R3").Value = 0.00584934642747072
$ws.Range("S3").Value = 0.00584934642747072
$ws.Range("U3").Value = 0.00584934642747072
$ws.Range("AK3").Value = 0.00584934642747072
$ws.Range("AM3").Value = 0.00584934642747072

$ws.Range("G4").Value = 0.005249676530928998
$ws.Range("I4").Value = 0.005249676530928998
$ws.Range("M4").Value = 0.005249676530928998
$ws.Range("O4").Value = 0.005249676530928998
$ws.Range("P4").Value = 0.005249676530928998
$ws.Range("R4").Value = 0.005249676530928998
$ws.Range("S4").Value = 0.005249676530928998
$ws.Range("U4").Value = 0.005249676530928998
$ws.Range("AK4").Value = 0.005249676530928998
$ws.Range("AM4").Value = 0.005249676530928998

$ws.Range("G5").Value = 0.003092023956230371
$ws.Range("I5").Value = 0.003092023956230371
$ws.Range("M5").Value = 0.003092023956230371
$ws.Range("O5").Value = 0.003092023956230371
$ws.Range("P5").Value = 0.003092023956230371
$ws.Range("R5").Value = 0.003092023956230371
$ws.Range("S5").Value = 0.003092023956230371
$ws.Range("U5").Value = 0.003092023956230371
$ws.Range("AK5").Value = 0.003092023956230371
$ws.Range("AM5").Value = 0.003092023956230371

$ws.Range("G6").Value = 0.005373854080857359
$ws.Range("I6").Value = 0.005373854080857359
$ws.Range("M6").Value = 0.005373854080857359
$ws.Range("O6").Value = 0.005373854080857359
$ws.Range("P6").Value = 0.005373854080857359
$ws.Range("R6").Value = 0.005373854080857359
$ws.Range("S6").Value = 0.005373854080857359
$ws.Range("U6").Value = 0.005373854080857359
$ws.Range("AK6").Value = 0.005373854080857359
$ws.Range("AM6").Value = 0.005373854080857359

$ws.Range("G7").Value = 0.003418627722410212
$ws.Range("I7").Value = 0.003418627722410212
$ws.Range("M7").Value = 0.003418627722410212
$ws.Range("O7").Value = 0.003418627722410212
$ws.Range("P7").Value = 0.003418627722410212
$ws.Range("R7").Value = 0.003418627722410212
$ws.Range("S7").Value = 0.003418627722410212
$ws.Range("U7").Value = 0.003418627722410212
$ws.Range("AK7").Value = 0.003418627722410212
$ws.Range("AM7").Value = 0.003418627722410212

$ws.Range("G8").Value = 0.001868579563410878
$ws.Range("I8").Value = 0.001868579563410878
$ws.Range("M8").Value = 0.001868579563410878
$ws.Range("O8").Value = 0.001868579563410878
$ws.Range("P8").Value = 0.001868579563410878
$ws.Range("R8").Value = 0.001868579563410878
$ws.Range("S8").Value = 0.001868579563410878
$ws.Range("U8").Value = 0.001868579563410878
$ws.Range("AK8").Value = 0.001868579563410878
$ws.Range("AM8").Value = 0.001868579563410878

$ws.Range("G9").Value = 0.001455408475076304
$ws.Range("I9").Value = 0.001455408475076304
$ws.Range("M9").Value = 0.001455408475076304
$ws.Range("O9").Value = 0.001455408475076304
$ws.Range("P9").Value = 0.001455408475076304
$ws.Range("R9").Value = 0.001455408475076304
$ws.Range("S9").Value = 0.001455408475076304
$ws.Range("U9").Value = 0.001455408475076304
$ws.Range("AK9").Value = 0.001455408475076304
$ws.Range("AM9").Value = 0.001455408475076304

$ws.Range("G10").Value = 0.01342174696481772
$ws.Range("I10").Value = 0.01342174696481772
$ws.Range("M10").Value = 0.01342174696481772
$ws.Range("O10").Value = 0.01342174696481772
$ws.Range("P10").Value = 0.01342174696481772
$ws.Range("R10").Value = 0.01342174696481772
$ws.Range("S10").Value = 0.01342174696481772
$ws.Range("U10").Value = 0.01342174696481772
$ws.Range("AK10").Value = 0.01342174696481772
$ws.Range("AM10").Value = 0.01342174696481772
